$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '67.317.02'
$ws.Range("E2").Value = '  -0.10%  '

$ws.Range("D3").Value = '3.487.25'
$ws.Range("E3").Value = '  -0.46%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.25%  '

$ws.Range("D7").Value = '3.486.80'
$ws.Range("E7").Value = '  -0.44%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("E9").Value = '  +0.90%  '

$ws.Range("E10").Value = '  +3.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.51'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.79%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.430'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '32.41'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.84%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000215'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.39%  '

$ws.Range("D15").Value = '4.079.91'
$ws.Range("E15").Value = '  -0.39%  '

$ws.Range("D16").Value = '3.491.50'
$ws.Range("E16").Value = '  -0.32%  '

$ws.Range("D17").Value = '67.290.21'
$ws.Range("E17").Value = '  -0.13%  '

$ws.Range("E18").Value = '  -0.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.53'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.70%  '

$ws.Range("E20").Value = '  +2.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '445.85'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.629'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.88%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.40'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("D25").Value = '3.634.27'
$ws.Range("E25").Value = '  -0.19%  '

$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000127'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.12%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.81'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.83%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.04'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.75%  '

$ws.Range("E30").Value = '  +0.31%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.63'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.47%  '

$ws.Range("E32").Value = '  +1.71%  '

$ws.Range("E33").Value = '  +0.10%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.66'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.15'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.03%  '

$ws.Range("E36").Value = '  +1.55%  '

$ws.Range("D37").Value = '3.481.85'
$ws.Range("E37").Value = '  -0.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.98'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.67%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.30'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.72%  '

$ws.Range("E41").Value = '  +0.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '174.55'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.75%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0894'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.83%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.46'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.90%  '

$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '29.83'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.73%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.875'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.31%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.74'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.04%  '

$ws.Range("E48").Value = '  +3.77%  '

$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.51'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.21%  '

$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.62'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.82%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.253'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.35%  '
